# "Correct error in BFoCPAbS"
#
# Root-cause edit: on the Rhodium sheet, the 2037 CCS-capacity figure (47)
# was being added to the wrong base case. It belongs to B35 (2037 capacity
# reported for the 2037 case, i.e. 28+47=75), not to B36 (which should just
# be the bare 19 figure).
#
# Every other change in this workbook (BAU Calculations rows 34/35/40-43,
# the various BFoCPAbS-* sheets, and the chart caches) is a pure downstream
# recalculation cascade from these two cells, so we only need to touch the
# Rhodium sheet and let Excel recalc the rest.

$wb = $excel.ActiveWorkbook

$rhodium = $wb.Worksheets.Item("Rhodium")

# B35: was a hard-coded 28 -> now a formula 28+47 (=75)
$rhodium.Range("B35").Formula = "=28+47"

# B36: was a formula 19+47 (=66) -> now just 19
$rhodium.Range("B36").Formula = "=19"

# Restore the window/selection state recorded for this sheet in the target
# workbook (topLeftCell A4 -> A25, active cell B35 -> B36).
$rhodium.Application.ActiveWindow.ScrollRow = 25
$rhodium.Range("B36").Select()

$wb.Application.Calculate()
